$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the crypto price (D) / 1h volume change (E) columns with the
# latest scrape. Plain assignment is used whenever possible; values that
# look like plain decimals (e.g. "20.01") would otherwise be auto-coerced
# to Number by Excel's input parser, so those are entered with a leading
# apostrophe to force Text, then ClearFormats() removes the quote-prefix
# style that entry method leaves behind -- keeping the cell unstyled, just
# like the rest of the sheet.

$ws.Range("D2").Value = '26.301.41'
$ws.Range("E2").Value = '  -4.43%  '
$ws.Range("D3").Value = '1.755.82'
$ws.Range("E3").Value = '  -4.17%  '
$ws.Range("E4").Value = '  +0.11%  '
$ws.Range("E5").Value = '  +0.09%  '
$ws.Range("D6").Value = '''302.98'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -3.08%  '
$ws.Range("D7").Value = '''0.4277'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  +0.51%  '
$ws.Range("D8").Value = '''0.3597'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("D9").Value = '''0.06981'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -3.90%  '
$ws.Range("D10").Value = '''0.8262'
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = '  -4.15%  '
$ws.Range("D11").Value = '''20.01'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  -2.88%  '
$ws.Range("D12").Value = '1.746.95'
$ws.Range("E12").Value = '  -3.54%  '
$ws.Range("D13").Value = '''5.188'
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = '  -3.68%  '
$ws.Range("D14").Value = '''6.313'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -2.87%  '
$ws.Range("D15").Value = '''0.06782'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.18%  '
$ws.Range("E16").Value = '  +0.28%  '
$ws.Range("D17").Value = '''78.68'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -2.17%  '
$ws.Range("D18").Value = '''0.000008623'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -2.98%  '
$ws.Range("D19").Value = '''1.005'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +0.39%  '
$ws.Range("E20").Value = '  -3.60%  '
$ws.Range("D21").Value = '26.268.14'
$ws.Range("E21").Value = '  -5.09%  '
$ws.Range("D22").Value = '''4.951'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -3.73%  '
$ws.Range("D23").Value = '''11.05'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.22%  '
$ws.Range("D24").Value = '1.969.08'
$ws.Range("E24").Value = '  -3.58%  '
$ws.Range("D25").Value = '''1.901'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -4.32%  '
$ws.Range("D26").Value = '''151.75'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -1.72%  '
$ws.Range("D27").Value = '''18.01'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -4.26%  '
$ws.Range("D28").Value = '''114.67'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  +0.31%  '
$ws.Range("D29").Value = '''4.991'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -2.08%  '
$ws.Range("D30").Value = '''1.643'
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = '  -9.44%  '
$ws.Range("D31").Value = '''0.08899'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  +0.60%  '
$ws.Range("D32").Value = '''0.7129'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -3.98%  '
$ws.Range("D33").Value = '''4.275'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  -5.53%  '
$ws.Range("D34").Value = '''1.087'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  -3.60%  '
$ws.Range("E35").Value = '  +0.04%  '
$ws.Range("E36").Value = '  -8.23%  '
$ws.Range("D37").Value = '''1.061'
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = '  -2.70%  '
$ws.Range("D38").Value = '''0.05066'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -4.34%  '
$ws.Range("D39").Value = '''0.01872'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -3.05%  '
$ws.Range("E40").Value = '  -3.75%  '
$ws.Range("D41").Value = '''0.4859'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -3.96%  '
$ws.Range("D42").Value = '''6.101'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  -5.69%  '
$ws.Range("D43").Value = '''2.459'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -12.06%  '
$ws.Range("E44").Value = '  -5.34%  '
$ws.Range("D45").Value = '''104.29'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("E46").Value = '  +0.17%  '
$ws.Range("D47").Value = '''9.918'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.72%  '
$ws.Range("D48").Value = '''0.06168'
$ws.Range("D48").ClearFormats()
$ws.Range("D49").Value = '''0.4439'
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = '  -4.86%  '
$ws.Range("D50").Value = '''1.558'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.14%  '
$ws.Range("D51").Value = '''1.702'
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = '  -0.54%  '
